$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '25.696.18'
$ws.Range("E2").Value = '  -3.62%  '

$ws.Range("D3").Value = '1.745.85'
$ws.Range("E3").Value = '  -5.61%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.001'
$ws.Range("D4").NumberFormat = "General"
$ws.Range("E4").Value = '  -0.01%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '236.84'
$ws.Range("D5").NumberFormat = "General"
$ws.Range("E5").Value = '  -9.82%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.000'
$ws.Range("D6").NumberFormat = "General"
$ws.Range("E6").Value = '  -0.09%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4930'
$ws.Range("D7").NumberFormat = "General"
$ws.Range("E7").Value = '  -8.31%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '41.59'
$ws.Range("D8").NumberFormat = "General"
$ws.Range("E8").Value = '  -7.72%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.2495'
$ws.Range("D9").NumberFormat = "General"
$ws.Range("E9").Value = '  -21.79%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.05996'
$ws.Range("D10").NumberFormat = "General"
$ws.Range("E10").Value = '  -14.47%  '

$ws.Range("D11").Value = '1.744.45'
$ws.Range("E11").Value = '  -5.75%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.06812'
$ws.Range("D12").NumberFormat = "General"
$ws.Range("E12").Value = '  -13.13%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '14.84'
$ws.Range("D13").NumberFormat = "General"
$ws.Range("E13").Value = '  -22.03%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '4.462'
$ws.Range("D14").NumberFormat = "General"
$ws.Range("E14").Value = '  -11.72%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '77.00'
$ws.Range("D15").NumberFormat = "General"
$ws.Range("E15").Value = '  -14.01%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.5637'
$ws.Range("D16").NumberFormat = "General"
$ws.Range("E16").Value = '  -27.41%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.9999'
$ws.Range("D17").NumberFormat = "General"
$ws.Range("E17").Value = '  -0.08%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '1.000'
$ws.Range("D18").NumberFormat = "General"
$ws.Range("E18").Value = '  -0.03%  '

$ws.Range("D19").Value = '25.735.34'
$ws.Range("E19").Value = '  -3.57%  '

$ws.Range("E20").Value = '  -19.52%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.000006535'
$ws.Range("D21").NumberFormat = "General"
$ws.Range("E21").Value = '  -18.44%  '

$ws.Range("D22").Value = '1.967.18'
$ws.Range("E22").Value = '  -5.79%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '3.988'
$ws.Range("D23").NumberFormat = "General"
$ws.Range("E23").Value = '  -14.21%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '5.016'
$ws.Range("D24").NumberFormat = "General"
$ws.Range("E24").Value = '  -17.02%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '7.858'
$ws.Range("D25").NumberFormat = "General"
$ws.Range("E25").Value = '  -16.53%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '136.60'
$ws.Range("D26").NumberFormat = "General"
$ws.Range("E26").Value = '  -4.55%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.807'
$ws.Range("D28").NumberFormat = "General"
$ws.Range("E28").Value = '  -18.53%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '14.68'
$ws.Range("D29").NumberFormat = "General"
$ws.Range("E29").Value = '  -14.25%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '101.58'
$ws.Range("D30").NumberFormat = "General"
$ws.Range("E30").Value = '  -9.03%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.748'
$ws.Range("D31").NumberFormat = "General"
$ws.Range("E31").Value = '  -13.34%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.08006'
$ws.Range("D32").NumberFormat = "General"
$ws.Range("E32").Value = '  -8.56%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.342'
$ws.Range("D33").NumberFormat = "General"
$ws.Range("E33").Value = '  -18.74%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.04399'
$ws.Range("D34").NumberFormat = "General"
$ws.Range("E34").Value = '  -9.85%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.9989'
$ws.Range("D35").NumberFormat = "General"
$ws.Range("E35").Value = '  -0.16%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.610'
$ws.Range("D36").NumberFormat = "General"
$ws.Range("E36").Value = '  -9.96%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.9811'
$ws.Range("D37").NumberFormat = "General"
$ws.Range("E37").Value = '  -14.22%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.6033'
$ws.Range("D38").NumberFormat = "General"
$ws.Range("E38").Value = '  -18.20%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.673'
$ws.Range("D39").NumberFormat = "General"
$ws.Range("E39").Value = '  -14.24%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.002'
$ws.Range("D40").NumberFormat = "General"
$ws.Range("E40").Value = '  -15.03%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.000'
$ws.Range("D41").NumberFormat = "General"
$ws.Range("E41").Value = '  -0.08%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '103.32'
$ws.Range("D42").NumberFormat = "General"
$ws.Range("E42").Value = '  -5.65%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.01498'
$ws.Range("D43").NumberFormat = "General"
$ws.Range("E43").Value = '  -14.46%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.7598'
$ws.Range("D44").NumberFormat = "General"
$ws.Range("E44").Value = '  -16.43%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '5.146'
$ws.Range("D45").NumberFormat = "General"
$ws.Range("E45").Value = '  -13.18%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.3714'
$ws.Range("D46").NumberFormat = "General"
$ws.Range("E46").Value = '  -23.29%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.05220'
$ws.Range("D47").NumberFormat = "General"
$ws.Range("E47").Value = '  -10.58%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.1071'
$ws.Range("D48").NumberFormat = "General"
$ws.Range("E48").Value = '  -14.40%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '30.08'
$ws.Range("D49").NumberFormat = "General"
$ws.Range("E49").Value = '  -14.48%  '

$ws.Range("B50").Value = 'Aptos'
$ws.Range("C50").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '5.871'
$ws.Range("D50").NumberFormat = "General"
$ws.Range("E50").Value = '  -24.26%  '

$ws.Range("B51").Value = 'Aave'
$ws.Range("C51").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '52.34'
$ws.Range("D51").NumberFormat = "General"
$ws.Range("E51").Value = '  -13.57%  '
